$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A (the original un-headered index column); remaining
# columns B:F shift left to become A:E.
$ws.Range("A1").EntireColumn.Delete()
